# Update "Case locations and outbreaks - Case alerts / public exposure sites" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old rows 16 and 17), shrinking the table from
# A1:D17 down to A1:D15.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(16).Delete()

# New data for rows 2-15 (Location, Site, Exposure period, Notes)
$data = @(
    @("Brandon Park",     "Kmart, Brandon Park Shopping Centre  Cnr Springvale Rd and Ferntree Gully Rd  Brandon Park, VIC 3170", "16:35-17:10  31/1/2021",  "Case attended venue"),
    @("Brighton",         "North Point Cafe  2B North Rd  Brighton, VIC 3186",                                                    "08:10-09:30  31/1/2021",  "Case dined outside and used bathroom"),
    @("Clayton South",    "Nakama Workshop  85 Main Rd  Clayton South, VIC 3169",                                                 "11:15-12:00  1/2/2021",   "Case attended venue"),
    @("Heatherton",       "Melbourne Golf Academy  385 Centre Dandenong Rd  Heatherton, VIC 3202",                                "17:19-18:35  1/2/2021",   "Case attended venue"),
    @("Keysborough",      "Aces Sporting Club (Driving Range)  Cnr Springvale Rd and Hutton Rd  Keysborough, VIC 3173",           "22:00-23:15  30/1/2021",  "Case attended venue"),
    @("Keysborough",      "Kmart, Parkmore Keysborough Shopping Centre  C/317 Cheltenham Rd  Keysborough, VIC 3173",              "16:00-17:00  31/1/2021",  "Case attended venue"),
    @("Melbourne",        "Exford Hotel  199 Russell St  Melbourne, VIC 3000",                                                    "23:00-23:35  29/1/2021",  "Case attended bottle shop"),
    @("Moorabbin Airport","Lululemon, DFO Moorabbin  Shop G-039/250 Centre Dandenong Rd  Moorabbin VIC 3194",                     "17:00-17:45  1/2/2021",   "Case attended venue"),
    @("Noble Park",       "Club Noble  46/56 Moodemere St  Noble Park VIC 3174",                                                  "14:36-15:30  30/01/2021", "Case attended venue"),
    @("Springvale",       "Bunnings Springvale  849 Princes Hwy  Springvale, VIC 3171",                                           "11:30-12:15  1/2/2021",   "Case attended venue"),
    @("Springvale",       "Coles Springvale  825 Dandenong Rd  Springvale, VIC 3171",                                             "17:00-18:00  31/1/2021",  "Case attended venue"),
    @("Springvale",       "Sharetea Springvale  27C Buckingham Ave  Springvale, VIC 3171",                                        "18:50-19:30  1/2/2021",   "Case attended venue"),
    @("Springvale",       "Woolworths Springvale  302 Springvale Rd  Springvale, VIC 3171",                                       "18:30-19:30  1/2/2021",   "Case attended venue"),
    @("West Melbourne",   "Kebab Kingz  438 Spencer St  West Melbourne, VIC 3003",                                                "23:24-00:15  29/1/2021",  "Case dined outside")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
